$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C) for rows 2-19 from 45177 (2023-09-08)
# to 45178 (2023-09-09), keeping the existing date formatting.
$ws.Range("C2:C19").Value = 45178
